# Refresh the letter-frequency table on Sheet1: recompute counts/percentages
# from a larger corpus, re-sort rows by descending frequency, and add a new
# row for the separately-tracked letter "ё" (so the table grows from 33 to
# 34 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row order (letters), counts and percentages, already sorted by
# descending frequency to match the refreshed corpus.
$letters = @("о","е","а","н","и","т","с","в","л","р","к","м","д","п","у","я","ь","ч","г","з","ы","б","ж","й","ш","х","ю","э","щ","ц","ф","ё","ъ")
$counts  = @(115667,91838,81731,68872,65574,64769,54653,48418,47693,40108,33391,32127,31566,27828,26970,24528,23157,19631,19295,17831,17825,17200,12116,10105,8395,7477,6147,3674,3024,2979,1877,836,307)
$pcts    = @(0.1125593489352468,0.08937056798840805,0.07953511500969727,0.06702160062825452,0.06381220872919564,0.06302883684358546,0.05318462566988028,0.04711714280431565,0.04641162154087791,0.03903040942615333,0.03249387656199975,0.03126383673167518,0.03071790924369094,0.02708033892268363,0.0262453909998842,0.02386900075807043,0.02253483572059023,0.01910356954833988,0.01877659693521563,0.01735193054945996,0.01734609175279703,0.01673788376707483,0.0117904767280162,0.009833506713156463,0.008169449664220535,0.007276113774791774,0.005981847181174941,0.003575289823269356,0.002942753518118272,0.002898962543146275,0.001826570222720899,0.0008135390017020093,0.0002987517625867426)

# The table grows by one row (33 -> 34 data rows). Give the new label cell
# the same look (bold, centered, bordered) as the rest of column A by
# copying the formatting down from the last existing row.
$ws.Cells.Item(33, 1).Copy() | Out-Null
$ws.Cells.Item(34, 1).PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt $letters.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $letters[$i]
    $ws.Cells.Item($row, 2).Value = $counts[$i]
    $ws.Cells.Item($row, 3).Value = $pcts[$i]
}
